$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A26").Value = "sg_rr_36_025 2023-12-13 16-41-08"
$ws.Range("B26").Value = 0.01
$ws.Range("C26").Value = 1000
$ws.Range("D26").Value = 5001
$ws.Range("E26").Value = 1530
$ws.Range("F26").Value = 1570
$ws.Range("G26").Value = 0.5
$ws.Range("H26").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I26").Value = 2.5
$ws.Range("J26").Value = 2.5459999999999998
$ws.Range("K26").Value = 0.12551152479809299
$ws.Range("L26").Value = "prominence set by looking at roughly biggest height span of noise bits that don't appear visually to contain resonance peaks"

$ws.Range("H27").Select()
